$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,45,46)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "No"
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 10).Value = 0
}
